# Populate the "Catcher's View" strikezone visual columns (Pitch / Choice / Result)
# for each recorded pitch, reorder the "Pitch Mix:" lists to lead with CH, and fill
# in a couple of previously-blank outcome fields (Result + Exit Velo + Launch Angle).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- At-bat 1 (rows 10-16), Pitch Mix row 17 ---
$ws.Range("F10").Value = "CH"
$ws.Range("G10").Value = "Take"
$ws.Range("H10").Value = "Strike"

$ws.Range("F11").Value = "CB"
$ws.Range("G11").Value = "Take"
$ws.Range("H11").Value = "Ball"

$ws.Range("F12").Value = "FB"
$ws.Range("G12").Value = "Swing"
$ws.Range("H12").Value = "Strike"
$ws.Range("M12").Value = $null

$ws.Range("F13").Value = "CB"
$ws.Range("G13").Value = "Take"
$ws.Range("H13").Value = "Ball"

$ws.Range("F14").Value = "CB"
$ws.Range("G14").Value = "Take"
$ws.Range("H14").Value = "Ball"

$ws.Range("F15").Value = "FB"
$ws.Range("G15").Value = "Swing"
$ws.Range("H15").Value = "Foul"
$ws.Range("M15").Value = "Walk"

$ws.Range("F16").Value = "FB"
$ws.Range("G16").Value = "Take"
$ws.Range("H16").Value = "Ball"

$ws.Range("J17").Value = "CH,CB,FB"

# --- At-bat 2 (rows 19-25), Pitch Mix row 26 ---
$ws.Range("F19").Value = "CB"
$ws.Range("G19").Value = "Swing"
$ws.Range("H19").Value = "Strike"
$ws.Range("M19").Value = "81.91 MPH"

$ws.Range("F20").Value = "CB"
$ws.Range("G20").Value = "Take"
$ws.Range("H20").Value = "Ball"

$ws.Range("F21").Value = "FB"
$ws.Range("G21").Value = "Swing"
$ws.Range("H21").Value = "In Play"
$ws.Range("M21").Value = "0.22°"

$ws.Range("J26").Value = "CH,CB,FB"

# --- At-bat 3 (rows 28-34), Pitch Mix row 35 ---
$ws.Range("F28").Value = "FB"
$ws.Range("G28").Value = "Take"
$ws.Range("H28").Value = "Strike"
$ws.Range("M28").Value = "85.8 MPH"

$ws.Range("F29").Value = "CB"
$ws.Range("G29").Value = "Take"
$ws.Range("H29").Value = "Ball"

$ws.Range("F30").Value = "FB"
$ws.Range("G30").Value = "Take"
$ws.Range("H30").Value = "Ball"
$ws.Range("M30").Value = "52.25°"

$ws.Range("F31").Value = "CH"
$ws.Range("G31").Value = "Swing"
$ws.Range("H31").Value = "In Play"

$ws.Range("J35").Value = "CH,CB,FB,SL"

# --- At-bat 4 (rows 37-43), Pitch Mix row 44 ---
$ws.Range("F37").Value = "CH"
$ws.Range("G37").Value = "Take"
$ws.Range("H37").Value = "Ball"
$ws.Range("M37").Value = "92.32 MPH"

$ws.Range("F38").Value = "CH"
$ws.Range("G38").Value = "Swing"
$ws.Range("H38").Value = "In Play"

$ws.Range("M39").Value = "49.5°"

$ws.Range("J44").Value = "CH,CB,FB,SL"

# --- At-bat 5 (rows 46-52), Pitch Mix row 53 ---
$ws.Range("F46").Value = "SL"
$ws.Range("G46").Value = "Swing"
$ws.Range("H46").Value = "In Play"
$ws.Range("M46").Value = "64.03 MPH"

$ws.Range("M48").Value = "1.21°"

$ws.Range("J53").Value = "CH,FB,SL"

# --- At-bat 9 (rows 61-67), Pitch Mix row 68 ---
$ws.Range("G61").Value = "Take"
$ws.Range("H61").Value = "Ball"
$ws.Range("M61").Value = "97.64 MPH"

$ws.Range("F62").Value = "CH"
$ws.Range("G62").Value = "Swing"
$ws.Range("H62").Value = "Strike"

$ws.Range("F63").Value = "CH"
$ws.Range("G63").Value = "Swing"
$ws.Range("H63").Value = "In Play"
$ws.Range("M63").Value = "-7.1°"

$ws.Range("J68").Value = "CH,FB,SL"
